# [Kadastro App] Kayıt silindi: 11308442
# Remove the record row for Kayıt No 11308442 from both the master
# "Kayitlar" sheet and the filtered "Merkez İlçe" sheet (the record's
# Birim is "Merkez İlçe", so it also appears there).

$wb = $excel.ActiveWorkbook

# --- Master sheet: "Kayitlar" ---
$wsMaster = $wb.Worksheets.Item("Kayitlar")
$wsMaster.Rows.Item(1187).Delete()

# --- Filtered sheet: "Merkez İlçe" ---
$wsMerkez = $wb.Worksheets.Item("Merkez İlçe")
$wsMerkez.Rows.Item(648).Delete()
